$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.169.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.66%  "
$ws.Range("D3").Value = "'1.577.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  -1.02%  "
$ws.Range("D5").Value = "'213.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.00%  "
$ws.Range("D8").Value = "'23.50"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.73%  "
$ws.Range("D9").Value = "'0.250"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.55%  "
$ws.Range("D10").Value = "'0.0600"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.05%  "
$ws.Range("D11").Value = "'0.0883"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.81%  "
$ws.Range("D12").Value = "'1.803.33"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("D13").Value = "'1.591.40"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.43%  "
$ws.Range("E14").Value = "  -0.65%  "
$ws.Range("D15").Value = "'0.525"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.14%  "
$ws.Range("D16").Value = "'28.129.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.54%  "
$ws.Range("D17").Value = "'63.68"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.17%  "
$ws.Range("D18").Value = "'230.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.69%  "
$ws.Range("E19").Value = "  +0.59%  "
$ws.Range("D20").Value = "'7.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.86%  "
$ws.Range("E21").Value = "  -0.97%  "
$ws.Range("E22").Value = "  -0.53%  "
$ws.Range("D23").Value = "'9.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.32%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").Value = "'152.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.21%  "
$ws.Range("D26").Value = "'15.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.21%  "
$ws.Range("E27").Value = "  -0.98%  "
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("E29").Value = "  -1.01%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").Value = "'0.0475"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.40%  "
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("E33").Value = "  -1.53%  "
$ws.Range("D34").Value = "'1.418.24"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.04%  "
$ws.Range("E35").Value = "  -1.29%  "
$ws.Range("D36").Value = "'1.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.34%  "
$ws.Range("E37").Value = "  -1.23%  "
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("E39").Value = "  +1.28%  "
$ws.Range("D40").Value = "'2.49"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.92%  "
$ws.Range("E42").Value = "  -1.10%  "
$ws.Range("E43").Value = "  -2.57%  "
$ws.Range("E44").Value = "  -2.71%  "
$ws.Range("D45").Value = "'1.82"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.03%  "
$ws.Range("D46").Value = "'63.93"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.20%  "
$ws.Range("D47").Value = "'1.716.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.63%  "
$ws.Range("D48").Value = "'87.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.44%  "
$ws.Range("D49").Value = "'0.0₆0108"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.77%  "
$ws.Range("E50").Value = "  +0.98%  "
$ws.Range("D51").Value = "'0.0944"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.65%  "
